$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "2014 Cohort" to "2014_Cohort"
$ws.Name = "2014_Cohort"

# Move the active selection to E19 (was C1:C1048576 / activeCell C1)
$ws.Range("E19").Select()
